$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "Meters": rename "Accu Chek Compact Plus" to "Accu Chek Compact Plus*"
$ws1.Range("D2").Value = "Accu Chek Compact Plus*"
# column D grew a bit wider once Excel re-measured the longer label
$ws1.Columns.Item(4).ColumnWidth = 22.83

# --- Sheet "VID AND PID": bump VID/PID key for Bayer row, add a second
#     Bayer device (Contour USB / Contour USB Next) on a new row
$ws2.Range("C4").Value = 6002
$ws2.Range("F4").Value = "Contour USB"
$ws2.Range("B5").Value = "1a79"
$ws2.Range("C5").Value = 7410
$ws2.Range("F5").Value = "Contour USB Next"

# --- Selections / active sheet: author left off on "VID AND PID" with G5
#     selected, having last looked at D9 on "Meters"
$ws1.Range("D9").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("G5").Select() | Out-Null
